$wb = $excel.ActiveWorkbook

# Update the metadata "Date" value on the Metadata sheet (row 8, col B)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-04-05T11:13:11-04:00"

# Remove the "Include from LOINC" sheet entirely
$loinc = $wb.Worksheets.Item("Include from LOINC")
$loinc.Delete()

# Keep the originally active sheet selected
$meta.Activate()
